# "Fruta / hortaliza, semanal"
#
# A new weekly price record is inserted at row 68 of the Cilantro sheet
# (Agrícola del Norte S.A. de Arica). All the existing records from row 68
# down get pushed one row lower (68->69, ..., 92->93), and the freshly
# inserted row 68 carries the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 68..92 down to 69..93, opening up a blank row 68.
$ws.Rows.Item(68).Insert()

# Populate the new row 68 with the new weekly record. Every column except
# the date (D) and the min/max/weighted-avg/unit prices (K, L, M, P) is the
# same constant used throughout this sheet (same market, same product).
$ws.Range("A68").Value = 1
$ws.Range("B68").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C68").Value = "Arica y Parinacota"
$ws.Range("D68").Value = 44875
$ws.Range("E68").Value = 15
$ws.Range("F68").Value = 100112040
$ws.Range("G68").Value = "Cilantro"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 300
$ws.Range("K68").Value = 500
$ws.Range("L68").Value = 600
$ws.Range("M68").Value = 550
$ws.Range("N68").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 275
$ws.Range("Q68").Value = 2
$ws.Range("R68").Value = "Hortaliza"
